# Payment Request Form template restructure:
#  - Move the "Payment Details" table (header + 2 detail rows + total row)
#    from rows 7-10 down to rows 11-14.
#  - Move the "Qualified Receiver Names / Date Payment Authorized /
#    Expense Authority Name / Account Coding" rows from rows 11-14 up to
#    rows 7-10.
#  - Add a new blank row 15.
#  - Invoice Date / Invoice Number rows (5-6) keep their text but move to
#    a plain (non-bold) sub-style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Break the merges that currently live on rows 11-14 - that block is
#    being replaced by the (currently unmerged) Payment Details table.
# ---------------------------------------------------------------------
$ws.Range("B11:E11").UnMerge()
$ws.Range("B12:E12").UnMerge()
$ws.Range("B13:E13").UnMerge()
$ws.Range("B14:E14").UnMerge()

# ---------------------------------------------------------------------
# 2. Rows 7-10 : Qualified Receiver Names / Date Payment Authorized /
#    Expense Authority Name / Account Coding (formerly rows 11-14).
# ---------------------------------------------------------------------
$ws.Range("A7").Value2 = "Qualified Receiver Names"
$ws.Range("A7").Font.Bold = $true
$ws.Range("B7").Value2 = "{d.qualified_receiver_name}"
$ws.Range("B7").Font.Bold = $false
$ws.Range("B7:E7").Merge()

$ws.Range("A8").Value2 = "Date Payment Authorized"
$ws.Range("A8").Font.Bold = $true
$ws.Range("B8").Value2 = "{d.date_payment_authorized}"
$ws.Range("B8").Font.Bold = $false
$ws.Range("B8:E8").Merge()
$ws.Rows(8).RowHeight = 15

$ws.Range("A9").Value2 = "Expense Authority Name"
$ws.Range("A9").Font.Bold = $true
$ws.Range("B9").Value2 = "{d.expense_authority_name}"
$ws.Range("B9").Font.Bold = $false
$ws.Range("B9:E9").Merge()

$ws.Range("A10").Value2 = "Account Coding"
$ws.Range("A10").Font.Bold = $true
$ws.Range("B10").Value2 = "{d.account_coding}"
$ws.Range("B10").Font.Bold = $false
$ws.Range("B10:E10").Merge()

# ---------------------------------------------------------------------
# 3. Rows 11-14 : Payment Details table (formerly rows 7-10).
# ---------------------------------------------------------------------
$ws.Range("A11").Value2 = "Payment Details"
$ws.Range("A11").Font.Bold = $true
$ws.Range("B11").Value2 = "Agreement Number"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value2 = "Unique ID"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value2 = "Amount"
$ws.Range("D11").Font.Bold = $true
$ws.Range("E11").Value2 = ""
$ws.Range("E11").Font.Bold = $true
$ws.Rows(11).RowHeight = $ws.StandardHeight
$ws.Rows(11).AutoFit()

$ws.Range("A12").Value2 = " "
$ws.Range("A12").Font.Bold = $true
$ws.Range("B12").Value2 = "{d.payment_details[i].agreement_number}"
$ws.Range("B12").Font.Bold = $false
$ws.Range("C12").Value2 = "{d.payment_details[i].unique_id}"
$ws.Range("C12").Font.Bold = $false
$ws.Range("D12").Value2 = "{d.payment_details[i].amount}"
$ws.Range("D12").Font.Bold = $false
$ws.Rows(12).AutoFit()

$ws.Range("A13").Value2 = " "
$ws.Range("A13").Font.Bold = $true
$ws.Range("B13").Value2 = "{d.payment_details[i+1].agreement_number}"
$ws.Range("B13").Font.Bold = $false
$ws.Range("C13").Value2 = "{d.payment_details[i+1].unique_id}"
$ws.Range("C13").Font.Bold = $false
$ws.Range("D13").Value2 = "{d.payment_details[i+1].amount}"
$ws.Range("D13").Font.Bold = $false
$ws.Rows(13).AutoFit()

$ws.Range("A14").Value2 = " "
$ws.Range("A14").Font.Bold = $true
$ws.Range("B14").Value2 = "Total Amount"
$ws.Range("B14").Font.Bold = $true
$ws.Range("C14").Value2 = " "
$ws.Range("C14").Font.Bold = $false
$ws.Range("D14").Value2 = "{d.total_payment}"
$ws.Range("D14").Font.Bold = $false
$ws.Rows(14).AutoFit()

# ---------------------------------------------------------------------
# 4. New trailing blank row 15.
# ---------------------------------------------------------------------
$ws.Range("A15").Font.Bold = $true
$ws.Range("D15").Font.Bold = $false

# ---------------------------------------------------------------------
# 5. Invoice Date / Invoice Number rows keep their values, just drop to a
#    non-highlighted sub style (still not bold).
# ---------------------------------------------------------------------
$ws.Range("A5").Value2 = "Invoice Date"
$ws.Range("B5").Value2 = "{d.invoice_date}"
$ws.Range("A6").Value2 = "Invoice Number"
$ws.Range("B6").Value2 = "{d.invoice_number}"

# ---------------------------------------------------------------------
# 6. View / selection bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("B15").Select()

$wb.Windows.Item(1).WindowState = -4143
